# IDAHO_2024.xlsx cleanup script
# 1) Rename header columns to short machine-friendly names
# 2) Title-case the "de"/"del"/"el"/"los"/"las"/"y" connector words (and every
#    other word) throughout the Estado/Municipio columns
# 3) Normalize the "TOTAL" grand-total label to "Total"
# 4) Drop the trailing metadata/footnote rows (794-798)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header row renames -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2) Title-case state (col A) and municipality (col B) text cells ------
for ($r = 2; $r -le 791; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Text
        if ($v.Length -gt 0) {
            $words = $v.Split(" ")
            $result = ""
            foreach ($w in $words) {
                if ($w.Length -gt 0) {
                    $first = $w.Substring(0,1).ToUpper()
                    $rest = $w.Substring(1)
                    $result = $result + $first + $rest + " "
                } else {
                    $result = $result + " "
                }
            }
            $result = $result.TrimEnd()
            $cell.Value = $result
        }
    }
}

# --- 3) Grand total label ---------------------------------------------------
$ws.Range("A792").Value = "Total"

# --- 4) Remove trailing metadata rows --------------------------------------
$ws.Range("A794:A798").EntireRow.Delete()
